$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emails = @{
    4 = "peter.griffin@gmail.com"
    5 = "ron.swanson@gmail.com"
    6 = "matt.murdock@gmail.com"
    7 = "davey.jones@gmail.com"
    8 = "nick.cage@gmail.com"
}

foreach ($row in 4..8) {
    $cell = $ws.Cells.Item($row, 4)
    $email = $emails[$row]
    $cell.Value = $email
    $ws.Hyperlinks.Add($cell, "mailto:$email")
    $cell.Style = "Hyperlink"
}

[void]$ws.Range("D8").Select()
